# Add an "ecoregion" lookup column (Piedmont / Coastal Plain) to the
# "location" sheet so it can be used to drive the new nitrogen/phosphorus
# calculation function.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("location")

# location_id -> ecoregion, keyed by row number (row 1 is the header).
# Row 16 (Luzon Branch / TLU01) intentionally has no ecoregion value yet.
$ecoregion = @{
  2  = "Piedmont"        # RCR01 - Rock Creek Upper
  3  = "Piedmont"        # RCR09 - Rock Creek Lower
  4  = "Piedmont"        # TBK01 - Battery Kemble Creek
  5  = "Piedmont"        # TBR01 - Broad Branch
  6  = "Piedmont"        # TDA01 - Dalecarlia Tributary
  7  = "Piedmont"        # TDO01 - Dumbarton Oaks
  8  = "Coastal Plain"   # TDU01 - Fort Dupont Tributary
  9  = "Piedmont"        # TFB01 - Foundry Branch
  10 = "Coastal Plain"   # TFC01 - Fort Chaplin Tributary
  11 = "Coastal Plain"   # TFD01 - Fort Davis Tributary
  12 = "Piedmont"        # TFE01 - Fenwick Branch
  13 = "Coastal Plain"   # TFS01 - Fort Stanton Tributary
  14 = "Coastal Plain"   # THR01 - Hickey Run
  15 = "Piedmont"        # TKV01 - Klingle Valley Run
  17 = "Piedmont"        # TMH01 - Melvin Hazen Valley Branch
  18 = "Coastal Plain"   # TNA01 - Nash Run
  19 = "Piedmont"        # TNS01 - Normanstone
  20 = "Coastal Plain"   # TOR01 - Oxon Run
  21 = "Coastal Plain"   # TPB01 - Pope Branch
  22 = "Piedmont"        # TPI01 - Pinehurst Branch
  23 = "Piedmont"        # TPO01 - Portal Branch
  24 = "Piedmont"        # TPY01 - Piney Branch
  25 = "Piedmont"        # TSO01 - Soapstone Creek
  26 = "Coastal Plain"   # TTX27 - Texas Avenue Tributary
  27 = "Coastal Plain"   # TWB01 - Watts Branch - Lower
  28 = "Coastal Plain"   # TWB05 - Watts Branch - Upper
}

for ($r = 2; $r -le 28; $r++) {
  if ($ecoregion.ContainsKey($r)) {
    $ws.Cells.Item($r, 3).Value = $ecoregion[$r]
  }
}

$ws.Range("C1").Value = "ecoregion"
$ws.Range("C1").Select() | Out-Null
